# Append a new data row (row 4) to the Sheet with another date/rate pair,
# matching the "Date" / "Rate" columns already present in rows 1-3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (not auto-converted to a date serial / number) so the
# new cells keep the same literal-text representation as the existing rows.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("B4").NumberFormat = "@"

$ws.Range("A4").Value = "2024-12-16"
$ws.Range("B4").Value = "7.284041"

# Restore the default cell style so the new cells don't carry a lingering
# "Text" number format and instead match the unstyled cells around them.
$ws.Range("A4:B4").Style = "Normal"
